$d = $word.ActiveDocument

$replacements = @(
    @{old = "734÷5="; new = "575÷9="},
    @{old = "352÷3="; new = "619÷2="},
    @{old = "275÷4="; new = "746÷2="},
    @{old = "946÷8="; new = "241÷5="},
    @{old = "216÷4="; new = "615÷7="},
    @{old = "917÷6="; new = "616÷3="},
    @{old = "912÷5="; new = "398÷8="},
    @{old = "719÷4="; new = "910÷3="},
    @{old = "659÷5="; new = "216÷9="},
    @{old = "366÷2="; new = "168÷6="},
    @{old = "350÷2="; new = "196÷6="},
    @{old = "857÷8="; new = "409÷8="},
    @{old = "247÷3="; new = "569÷3="},
    @{old = "909÷4="; new = "461÷3="},
    @{old = "443÷9="; new = "659÷9="},
    @{old = "250÷4="; new = "235÷2="},
    @{old = "437÷7="; new = "313÷6="},
    @{old = "106÷4="; new = "408÷2="},
    @{old = "376÷3="; new = "847÷5="},
    @{old = "797÷9="; new = "591÷5="},
    @{old = "430÷8="; new = "710÷7="},
    @{old = "266÷9="; new = "927÷7="},
    @{old = "650÷4="; new = "865÷9="},
    @{old = "495÷6="; new = "184÷9="},
    @{old = "678÷2="; new = "122÷8="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
